$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Simple, globally-unique text replacements (Find/Replace on the
#    whole document Range is safe for these because each old string
#    occurs exactly once in the document).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("SoT Improved from Amber/Green to Green", $true, $false, $false, $false, $false, $true, 1, $false, "SoT Improved from Green to Amber/Green", 2)
$d.Content.Find.Execute("A13 Improved from Amber to Amber/Green", $true, $false, $false, $false, $false, $true, 1, $false, "A13 Improved from Amber/Green to Amber", 2)
$d.Content.Find.Execute("Columbia Improved from Amber/Green to Green", $true, $false, $false, $false, $false, $true, 1, $false, "F9 Improved from Amber/Green to Amber", 2)
$d.Content.Find.Execute("3 project(s) in total have improved", $true, $false, $false, $false, $false, $true, 1, $false, "4 project(s) in total have improved", 2)
$d.Content.Find.Execute("1 project(s) in total have decreased", $true, $false, $false, $false, $false, $true, 1, $false, "0 project(s) in total have decreased", 2)
$d.Content.Find.Execute("Columbia Improved from Amber to Green", $true, $false, $false, $false, $false, $true, 1, $false, "Columbia Improved from Green to Amber", 2)
$d.Content.Find.Execute("SoT Improved from Amber to Green", $true, $false, $false, $false, $false, $true, 1, $false, "SoT Improved from Green to Amber", 2)
$d.Content.Find.Execute("1 project(s) in total are missing a rating", $true, $false, $false, $false, $false, $true, 1, $false, "0 project(s) in total are missing a rating", 2)

# ---------------------------------------------------------------------
# 2) Paragraph-index-based edits. Walk from the bottom of the document
#    upward so earlier paragraph indices stay valid while later ones
#    are inserted/removed.
# ---------------------------------------------------------------------

# --- "SRO Benefits RAG" block (3rd block) ---------------------------
# Paragraph 39: "A11 Has not provided a rating" -> delete entirely
$p39 = $d.Paragraphs.Item(39)
$p39.Range.Delete()

# Paragraph 33: "1 project(s) in total have improved" -> "4 project(s) ..."
# (the earlier global Find/Replace above skipped this one deliberately,
#  because the same old text also appears, unchanged, in the "SRO
#  Finance confidence" block further up the document)
$p33 = $d.Paragraphs.Item(33)
$p33.Range.Text = "4 project(s) in total have improved"

# Paragraph 32: "SoT Improved from Amber to Green" -> "SoT Improved from Green to Amber"
# then insert three new improvement lines right after it.
$p32 = $d.Paragraphs.Item(32)
$p32.Range.Text = "SoT Improved from Green to Amber"
$p32.Range.InsertParagraphAfter()

$p33b = $d.Paragraphs.Item(33)
$p33b.Range.Text = "A13 Improved from Green to Amber"
$p33b.Range.InsertParagraphAfter()

$p34b = $d.Paragraphs.Item(34)
$p34b.Range.Text = "F9 Improved from Green to Amber"
$p34b.Range.InsertParagraphAfter()

$p35b = $d.Paragraphs.Item(35)
$p35b.Range.Text = "Columbia Improved from Green to Amber"

# --- "Departmental DCA" block (1st block) ----------------------------
# Paragraph 11: "A11 Worsened from Amber to Amber/Red" -> delete entirely
$p11 = $d.Paragraphs.Item(11)
$p11.Range.Delete()

# Paragraph 7: "Columbia Improved from Amber/Green to Green" was already
# replaced above (now "F9 Improved from Amber/Green to Amber"). Insert
# the new "Columbia Improved from Green to Amber/Green" line after it.
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "Columbia Improved from Green to Amber/Green"
